$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.882461666666666
$ws.Range("H2").Value = 8.647385
$ws.Range("I2").Value = 0.513165610165437
$ws.Range("J2").Value = 0.513165610165437
$ws.Range("M2").Value = 0.8658956666666665
$ws.Range("N2").Value = 2.597687
$ws.Range("O2").Value = 0.281473303461206
$ws.Range("P2").Value = 0.281473303461206
$ws.Range("Q2").Value = 2.495911066499444
$ws.Range("R2").Value = 22.463199598495
$ws.Range("S2").Value = 0.144442419515951
$ws.Range("T2").Value = 0.144442419515951
$ws.Range("G3").Value = 2.882461666666666
$ws.Range("H3").Value = 8.647385
$ws.Range("I3").Value = 0.513165610165437
$ws.Range("J3").Value = 0.513165610165437
$ws.Range("O3").Value = 0.3999629207966763
$ws.Range("P3").Value = 0.3999629207966763
$ws.Range("Q3").Value = 3.546595246975
$ws.Range("R3").Value = 31.919357222775
$ws.Range("S3").Value = 0.2052472162941767
$ws.Range("T3").Value = 0.2052472162941767
$ws.Range("G4").Value = 2.882461666666666
$ws.Range("H4").Value = 8.647385
$ws.Range("I4").Value = 0.513165610165437
$ws.Range("J4").Value = 0.513165610165437
$ws.Range("M4").Value = 0.4220566666666667
$ws.Range("N4").Value = 1.26617
$ws.Range("O4").Value = 0.1371963029585455
$ws.Range("P4").Value = 0.1371963029585455
$ws.Range("Q4").Value = 1.216562162827778
$ws.Range("R4").Value = 10.94905946545
$ws.Range("S4").Value = 0.07040442452016417
$ws.Range("T4").Value = 0.07040442452016418
$ws.Range("G5").Value = 2.882461666666666
$ws.Range("H5").Value = 8.647385
$ws.Range("I5").Value = 0.513165610165437
$ws.Range("J5").Value = 0.513165610165437
$ws.Range("M5").Value = 0.5579403333333333
$ws.Range("N5").Value = 1.673821
$ws.Range("O5").Value = 0.1813674727835722
$ws.Range("P5").Value = 0.1813674727835722
$ws.Range("Q5").Value = 1.608241623120555
$ws.Range("R5").Value = 14.474174608085
$ws.Range("S5").Value = 0.09307154983514511
$ws.Range("T5").Value = 0.09307154983514511
$ws.Range("I6").Value = 0.2254722714492577
$ws.Range("J6").Value = 0.2254722714492577
$ws.Range("M6").Value = 0.8658956666666665
$ws.Range("N6").Value = 2.597687
$ws.Range("O6").Value = 0.281473303461206
$ws.Range("P6").Value = 0.281473303461206
$ws.Range("Q6").Value = 1.096641564343222
$ws.Range("R6").Value = 9.869774079088998
$ws.Range("S6").Value = 0.06346442508372431
$ws.Range("T6").Value = 0.06346442508372432
$ws.Range("I7").Value = 0.2254722714492577
$ws.Range("J7").Value = 0.2254722714492577
$ws.Range("O7").Value = 0.3999629207966763
$ws.Range("P7").Value = 0.3999629207966763
$ws.Range("S7").Value = 0.09018054824750615
$ws.Range("T7").Value = 0.09018054824750615
$ws.Range("I8").Value = 0.2254722714492577
$ws.Range("J8").Value = 0.2254722714492577
$ws.Range("M8").Value = 0.4220566666666667
$ws.Range("N8").Value = 1.26617
$ws.Range("O8").Value = 0.1371963029585455
$ws.Range("P8").Value = 0.1371963029585455
$ws.Range("Q8").Value = 0.5345273119988889
$ws.Range("R8").Value = 4.81074580799
$ws.Range("S8").Value = 0.03093396206250377
$ws.Range("T8").Value = 0.03093396206250378
$ws.Range("I9").Value = 0.2254722714492577
$ws.Range("J9").Value = 0.2254722714492577
$ws.Range("M9").Value = 0.5579403333333333
$ws.Range("N9").Value = 1.673821
$ws.Range("O9").Value = 0.1813674727835722
$ws.Range("P9").Value = 0.1813674727835722
$ws.Range("Q9").Value = 0.7066215752207776
$ws.Range("R9").Value = 6.359594176986999
$ws.Range("S9").Value = 0.04089333605552344
$ws.Range("T9").Value = 0.04089333605552344
$ws.Range("G10").Value = 0.8351260000000001
$ws.Range("H10").Value = 2.505378
$ws.Range("I10").Value = 0.1486777598158359
$ws.Range("J10").Value = 0.1486777598158359
$ws.Range("M10").Value = 0.8658956666666665
$ws.Range("N10").Value = 2.597687
$ws.Range("O10").Value = 0.281473303461206
$ws.Range("P10").Value = 0.281473303461206
$ws.Range("Q10").Value = 0.7231319845206666
$ws.Range("R10").Value = 6.508187860686
$ws.Range("S10").Value = 0.04184882020657509
$ws.Range("T10").Value = 0.04184882020657509
$ws.Range("G11").Value = 0.8351260000000001
$ws.Range("H11").Value = 2.505378
$ws.Range("I11").Value = 0.1486777598158359
$ws.Range("J11").Value = 0.1486777598158359
$ws.Range("O11").Value = 0.3999629207966763
$ws.Range("P11").Value = 0.3999629207966763
$ws.Range("Q11").Value = 1.02754320603
$ws.Range("R11").Value = 9.24788885427
$ws.Range("S11").Value = 0.05946559107344845
$ws.Range("T11").Value = 0.05946559107344845
$ws.Range("G12").Value = 0.8351260000000001
$ws.Range("H12").Value = 2.505378
$ws.Range("I12").Value = 0.1486777598158359
$ws.Range("J12").Value = 0.1486777598158359
$ws.Range("M12").Value = 0.4220566666666667
$ws.Range("N12").Value = 1.26617
$ws.Range("O12").Value = 0.1371963029585455
$ws.Range("P12").Value = 0.1371963029585455
$ws.Range("Q12").Value = 0.3524704958066667
$ws.Range("R12").Value = 3.17223446226
$ws.Range("S12").Value = 0.02039803897889129
$ws.Range("T12").Value = 0.0203980389788913
$ws.Range("G13").Value = 0.8351260000000001
$ws.Range("H13").Value = 2.505378
$ws.Range("I13").Value = 0.1486777598158359
$ws.Range("J13").Value = 0.1486777598158359
$ws.Range("M13").Value = 0.5579403333333333
$ws.Range("N13").Value = 1.673821
$ws.Range("O13").Value = 0.1813674727835722
$ws.Range("P13").Value = 0.1813674727835722
$ws.Range("Q13").Value = 0.4659504788153334
$ws.Range("R13").Value = 4.193554309338
$ws.Range("S13").Value = 0.02696530955692111
$ws.Range("T13").Value = 0.02696530955692111
$ws.Range("G14").Value = 0.6329503333333334
$ws.Range("H14").Value = 1.898851
$ws.Range("I14").Value = 0.1126843585694693
$ws.Range("J14").Value = 0.1126843585694693
$ws.Range("M14").Value = 0.8658956666666665
$ws.Range("N14").Value = 2.597687
$ws.Range("O14").Value = 0.281473303461206
$ws.Range("P14").Value = 0.281473303461206
$ws.Range("Q14").Value = 0.5480689508485556
$ws.Range("R14").Value = 4.932620557637
$ws.Range("S14").Value = 0.03171763865495558
$ws.Range("T14").Value = 0.03171763865495558
$ws.Range("G15").Value = 0.6329503333333334
$ws.Range("H15").Value = 1.898851
$ws.Range("I15").Value = 0.1126843585694693
$ws.Range("J15").Value = 0.1126843585694693
$ws.Range("O15").Value = 0.3999629207966763
$ws.Range("P15").Value = 0.3999629207966763
$ws.Range("Q15").Value = 0.7787852548850001
$ws.Range("R15").Value = 7.009067293965
$ws.Range("S15").Value = 0.04506956518154492
$ws.Range("T15").Value = 0.04506956518154492
$ws.Range("G16").Value = 0.6329503333333334
$ws.Range("H16").Value = 1.898851
$ws.Range("I16").Value = 0.1126843585694693
$ws.Range("J16").Value = 0.1126843585694693
$ws.Range("M16").Value = 0.4220566666666667
$ws.Range("N16").Value = 1.26617
$ws.Range("O16").Value = 0.1371963029585455
$ws.Range("P16").Value = 0.1371963029585455
$ws.Range("Q16").Value = 0.2671409078522223
$ws.Range("R16").Value = 2.40426817067
$ws.Range("S16").Value = 0.01545987739698629
$ws.Range("T16").Value = 0.01545987739698629
$ws.Range("G17").Value = 0.6329503333333334
$ws.Range("H17").Value = 1.898851
$ws.Range("I17").Value = 0.1126843585694693
$ws.Range("J17").Value = 0.1126843585694693
$ws.Range("M17").Value = 0.5579403333333333
$ws.Range("N17").Value = 1.673821
$ws.Range("O17").Value = 0.1813674727835722
$ws.Range("P17").Value = 0.1813674727835722
$ws.Range("Q17").Value = 0.3531485199634444
$ws.Range("R17").Value = 3.178336679671
$ws.Range("S17").Value = 0.02043727733598252
$ws.Range("T17").Value = 0.02043727733598251
